$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format price cells that look like plain numbers so Excel
# keeps them as text (matches the source data which stores prices
# as text, including multi-dot big numbers like "64.910.48").
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D47",
    "D50"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '64.910.48'
$ws.Range('E2').Value = '  +4.51%  '

$ws.Range('D3').Value = '3.106.34'
$ws.Range('E3').Value = '  +2.83%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = '559.58'
$ws.Range('E5').Value = '  +2.83%  '

$ws.Range('D6').Value = '145.20'
$ws.Range('E6').Value = '  +8.17%  '

$ws.Range('E7').Value = '  -0.09%  '

$ws.Range('D8').Value = '3.104.66'
$ws.Range('E8').Value = '  +3.00%  '

$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  +1.78%  '

$ws.Range('D10').Value = '7.15'
$ws.Range('E10').Value = '  +13.43%  '

$ws.Range('D11').Value = '0.153'
$ws.Range('E11').Value = '  +3.86%  '

$ws.Range('D12').Value = '0.468'
$ws.Range('E12').Value = '  +4.41%  '

$ws.Range('D13').Value = '0.0000231'
$ws.Range('E13').Value = '  +4.21%  '

$ws.Range('D14').Value = '35.50'
$ws.Range('E14').Value = '  +2.03%  '

$ws.Range('D15').Value = '3.596.95'
$ws.Range('E15').Value = '  +2.65%  '

$ws.Range('D16').Value = '64.867.80'
$ws.Range('E16').Value = '  +4.52%  '

$ws.Range('D17').Value = '3.097.12'
$ws.Range('E17').Value = '  +2.88%  '

$ws.Range('E18').Value = '  -0.54%  '

$ws.Range('D19').Value = '6.84'
$ws.Range('E19').Value = '  +2.69%  '

$ws.Range('D20').Value = '484.29'
$ws.Range('E20').Value = '  +0.08%  '

$ws.Range('D21').Value = '13.84'
$ws.Range('E21').Value = '  +4.11%  '

$ws.Range('D22').Value = '0.681'
$ws.Range('E22').Value = '  +1.17%  '

$ws.Range('D23').Value = '7.55'
$ws.Range('E23').Value = '  +7.56%  '

$ws.Range('D24').Value = '13.45'
$ws.Range('E24').Value = '  +11.68%  '

$ws.Range('D25').Value = '81.23'
$ws.Range('E25').Value = '  -1.10%  '

$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.04%  '

$ws.Range('D27').Value = '2.79'
$ws.Range('E27').Value = '  +3.27%  '

$ws.Range('D28').Value = '8.24'
$ws.Range('E28').Value = '  +6.16%  '

$ws.Range('E29').Value = '  +7.57%  '

$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  +0.15%  '

$ws.Range('D31').Value = '26.19'
$ws.Range('E31').Value = '  +1.48%  '

$ws.Range('E32').Value = '  +1.42%  '

$ws.Range('D33').Value = '2.47'
$ws.Range('E33').Value = '  +4.54%  '

$ws.Range('D34').Value = '5.72'
$ws.Range('E34').Value = '  +0.74%  '

$ws.Range('D35').Value = '6.25'
$ws.Range('E35').Value = '  +6.11%  '

$ws.Range('D36').Value = '55.02'
$ws.Range('E36').Value = '  -0.70%  '

$ws.Range('D37').Value = '469.43'
$ws.Range('E37').Value = '  +4.14%  '

$ws.Range('D38').Value = '0.0410'
$ws.Range('E38').Value = '  +6.83%  '

$ws.Range('D39').Value = '0.0830'
$ws.Range('E39').Value = '  +3.86%  '

$ws.Range('D40').Value = '2.90'
$ws.Range('E40').Value = '  +18.33%  '

$ws.Range('D41').Value = '3.025.49'
$ws.Range('E41').Value = '  -4.29%  '

$ws.Range('D42').Value = '8.28'
$ws.Range('E42').Value = '  +1.97%  '

$ws.Range('E43').Value = '  -1.42%  '

$ws.Range('D44').Value = '28.45'
$ws.Range('E44').Value = '  +7.65%  '

$ws.Range('D45').Value = '0.259'
$ws.Range('E45').Value = '  +6.16%  '

$ws.Range('E46').Value = '  +0.03%  '

$ws.Range('D47').Value = '2.12'
$ws.Range('E47').Value = '  +7.79%  '

$ws.Range('E48').Value = '  +3.76%  '

$ws.Range('E49').Value = '  +5.40%  '

$ws.Range('D50').Value = '118.00'
$ws.Range('E50').Value = '  +1.44%  '

$ws.Range('E51').Value = '  +2.14%  '

# Revert the temporary Text number format back to the default
# "Normal" style so we do not leave stray formatting behind.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
